$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.275.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.01%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.668.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.36%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'219.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.65%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5285"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.15%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2653"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.59%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.24%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07834"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.526"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.56%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.684.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.88%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.896.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.30%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.5603"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0₅8118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.80%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'65.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.292.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.17%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.727"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.53%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'200.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.29%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.80%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.068"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.16%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.86%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1217"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.24%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.239"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'16.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.525"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.78%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05910"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.94%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.517"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.85%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.327"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.599"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.9622"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.819"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.428"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5803"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.28%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.075.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.70%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8576"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'102.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'1.805.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.11%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'58.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.28%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.11%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'8.096"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.14%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₈103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.17%  "
$ws.Range("E50").Style = "Normal"
